$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 48113730
$ws.Range("I132").Value = 48113730
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 144341190
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -144338660
$ws.Range("N132").ClearContents()

$ws.Range("H135").Value = 3537.5
$ws.Range("I135").Value = 3570.5278
$ws.Range("J135").Value = 3339.3333
$ws.Range("K135").Value = 32134.7502
$ws.Range("L135").Value = 30053.9997
$ws.Range("M135").Value = -29599.7502
$ws.Range("N135").Value = -35123.9997

$ws.Range("H137").Value = 2452383.2
$ws.Range("I137").Value = 849.381
$ws.Range("J137").Value = 5026494
$ws.Range("K137").Value = 2548.143
$ws.Range("L137").Value = 15079482
$ws.Range("M137").Value = 1.856999999999971
$ws.Range("N137").Value = -15084582

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6540.2964
$ws.Range("I61").Value = 8543.111000000001
$ws.Range("J61").Value = 2534.6667
$ws.Range("K61").Value = 8543.111000000001
$ws.Range("L61").Value = 2534.6667
$ws.Range("M61").Value = -8331.111000000001
$ws.Range("N61").Value = -2958.6667

$ws.Range("H74").Value = 6091.3184
$ws.Range("I74").Value = 1233.4615
$ws.Range("J74").Value = 13108.223
$ws.Range("K74").Value = 1233.4615
$ws.Range("L74").Value = 13108.223
$ws.Range("M74").Value = -359.4614999999999
$ws.Range("N74").Value = -14856.223

$ws.Range("H77").Value = 6091.3184
$ws.Range("I77").Value = 1233.4615
$ws.Range("J77").Value = 13108.223
$ws.Range("K77").Value = 6167.307499999999
$ws.Range("L77").Value = 65541.11500000001
$ws.Range("M77").Value = -1799.307499999999
$ws.Range("N77").Value = -74277.11500000001

$ws.Range("H102").Value = 166668240
$ws.Range("I102").Value = 200001300
$ws.Range("K102").Value = 200001300
$ws.Range("M102").Value = -199999678

$ws.Range("H122").Value = 52175268
$ws.Range("I122").Value = 60001228
$ws.Range("J122").Value = 2204.3333
$ws.Range("K122").Value = 180003684
$ws.Range("L122").Value = 6612.999899999999
$ws.Range("M122").Value = -180001234
$ws.Range("N122").Value = -11512.9999

$ws.Range("H136").Value = 6540.2964
$ws.Range("I136").Value = 8543.111000000001
$ws.Range("J136").Value = 2534.6667
$ws.Range("K136").Value = 25629.333
$ws.Range("L136").Value = 7604.000100000001
$ws.Range("M136").Value = -23079.333
$ws.Range("N136").Value = -12704.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 25885.666
$ws.Range("J103").Value = 25885.666
$ws.Range("L103").Value = 25885.666
$ws.Range("N103").Value = -28229.666

$ws.Range("H134").Value = 9821200
$ws.Range("I134").Value = 15176684
$ws.Range("J134").Value = 2810.8333
$ws.Range("K134").Value = 45530052
$ws.Range("L134").Value = 8432.499899999999
$ws.Range("M134").Value = -45527517
$ws.Range("N134").Value = -13502.4999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6730.7847
$ws.Range("I31").Value = 3614.1052
$ws.Range("J31").Value = 9619.415000000001
$ws.Range("K31").Value = 3614.1052
$ws.Range("L31").Value = 9619.415000000001
$ws.Range("M31").Value = -3319.1052
$ws.Range("N31").Value = -10209.415

$ws.Range("H34").Value = 6730.7847
$ws.Range("I34").Value = 3614.1052
$ws.Range("J34").Value = 9619.415000000001
$ws.Range("K34").Value = 3614.1052
$ws.Range("L34").Value = 9619.415000000001
$ws.Range("M34").Value = -3412.1052
$ws.Range("N34").Value = -10023.415

$ws.Range("H58").Value = 3768995
$ws.Range("I58").Value = 7937290.5
$ws.Range("K58").Value = 7937290.5
$ws.Range("M58").Value = -7937087.5

$ws.Range("H99").Value = 2461.9524
$ws.Range("I99").Value = 2034.3077
$ws.Range("J99").Value = 3156.875
$ws.Range("K99").Value = 2034.3077
$ws.Range("L99").Value = 3156.875
$ws.Range("M99").Value = -536.3077000000001
$ws.Range("N99").Value = -6152.875

$ws.Range("H126").Value = 2461.9524
$ws.Range("I126").Value = 2034.3077
$ws.Range("J126").Value = 3156.875
$ws.Range("K126").Value = 6102.9231
$ws.Range("L126").Value = 9470.625
$ws.Range("M126").Value = -3632.9231
$ws.Range("N126").Value = -14410.625

$ws.Range("H132").Value = 5955924.5
$ws.Range("I132").Value = 9010447
$ws.Range("J132").Value = 7643.579
$ws.Range("K132").Value = 27031341
$ws.Range("L132").Value = 22930.737
$ws.Range("M132").Value = -27028811
$ws.Range("N132").Value = -27990.737

$ws.Range("H136").Value = 3768995
$ws.Range("I136").Value = 7937290.5
$ws.Range("K136").Value = 23811871.5
$ws.Range("M136").Value = -23809321.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 217390.55
$ws.Range("I68").Value = 580
$ws.Range("J68").Value = 486938.78
$ws.Range("K68").Value = 1740
$ws.Range("L68").Value = 1460816.34
$ws.Range("M68").Value = -929
$ws.Range("N68").Value = -1462438.34

$ws.Range("H71").Value = 217390.55
$ws.Range("I71").Value = 580
$ws.Range("J71").Value = 486938.78
$ws.Range("K71").Value = 5220
$ws.Range("L71").Value = 4382449.02
$ws.Range("M71").Value = -1164
$ws.Range("N71").Value = -4390561.02

$ws.Range("H107").Value = 838.4286
$ws.Range("I107").Value = 392.75
$ws.Range("J107").Value = 3512.5
$ws.Range("K107").Value = 1178.25
$ws.Range("L107").Value = 10537.5
$ws.Range("M107").Value = 741.75
$ws.Range("N107").Value = -14377.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 100000990
$ws.Range("I97").Value = 71429490
$ws.Range("J97").Value = 166667820
$ws.Range("K97").Value = 71429490
$ws.Range("L97").Value = 166667820
$ws.Range("M97").Value = -71428994
$ws.Range("N97").Value = -166668812

$ws.Range("H102").Value = 33334926
$ws.Range("I102").Value = 62501852
$ws.Range("J102").Value = 1294.8572
$ws.Range("K102").Value = 62501852
$ws.Range("L102").Value = 1294.8572
$ws.Range("M102").Value = -62500230
$ws.Range("N102").Value = -4538.8572

$ws.Range("H122").Value = 58824708
$ws.Range("I122").Value = 76923990
$ws.Range("J122").Value = 2025
$ws.Range("K122").Value = 230771970
$ws.Range("L122").Value = 6075
$ws.Range("M122").Value = -230769520
$ws.Range("N122").Value = -10975

$ws.Range("H126").Value = 2099.3125
$ws.Range("I126").Value = 1509.375
$ws.Range("J126").Value = 2689.25
$ws.Range("K126").Value = 4528.125
$ws.Range("L126").Value = 8067.75
$ws.Range("M126").Value = -2058.125
$ws.Range("N126").Value = -13007.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3211885.8
$ws.Range("I7").Value = 4011.75
$ws.Range("J7").Value = 9627634
$ws.Range("K7").Value = 4011.75
$ws.Range("L7").Value = 9627634
$ws.Range("M7").Value = -3899.75
$ws.Range("N7").Value = -9627858

$ws.Range("H40").Value = 3399.5
$ws.Range("I40").Value = 3366
$ws.Range("K40").Value = 3366
$ws.Range("M40").Value = -3230

$ws.Range("H100").Value = 1986.1428
$ws.Range("I100").Value = 1443.7142
$ws.Range("J100").Value = 2528.5715
$ws.Range("K100").Value = 1443.7142
$ws.Range("L100").Value = 2528.5715
$ws.Range("M100").Value = -902.7141999999999
$ws.Range("N100").Value = -3610.5715

$ws.Range("H122").Value = 31253850
$ws.Range("I122").Value = 1300
$ws.Range("J122").Value = 35718500
$ws.Range("K122").Value = 3900
$ws.Range("L122").Value = 107155500
$ws.Range("M122").Value = -1450
$ws.Range("N122").Value = -107160400

$ws.Range("H123").Value = 29958.215
$ws.Range("J123").Value = 29958.215
$ws.Range("L123").Value = 29958.215
$ws.Range("N123").Value = -39758.215

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H126").Value = 3211885.8
$ws.Range("I126").Value = 4011.75
$ws.Range("J126").Value = 9627634
$ws.Range("K126").Value = 12035.25
$ws.Range("L126").Value = 28882902
$ws.Range("M126").Value = -9565.25
$ws.Range("N126").Value = -28887842

$ws.Range("H132").Value = 2849069.5
$ws.Range("I132").Value = 3088543.8
$ws.Range("J132").Value = 2395328.8
$ws.Range("K132").Value = 9265631.399999999
$ws.Range("L132").Value = 7185986.399999999
$ws.Range("M132").Value = -9263101.399999999
$ws.Range("N132").Value = -7191046.399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 71428960
$ws.Range("I107").Value = 520.5
$ws.Range("J107").Value = 166666860
$ws.Range("K107").Value = 1561.5
$ws.Range("L107").Value = 500000580
$ws.Range("M107").Value = 358.5
$ws.Range("N107").Value = -500004420

$ws.Range("H122").Value = 9887.362999999999
$ws.Range("I122").Value = 13221.375
$ws.Range("J122").Value = 996.6667
$ws.Range("K122").Value = 39664.125
$ws.Range("L122").Value = 2990.0001
$ws.Range("M122").Value = -37214.125
$ws.Range("N122").Value = -7890.0001

$ws.Range("H136").Value = 13342980
$ws.Range("I136").Value = 6369094
$ws.Range("J136").Value = 33334788
$ws.Range("K136").Value = 19107282
$ws.Range("L136").Value = 100004364
$ws.Range("M136").Value = -19104732
$ws.Range("N136").Value = -100009464
